$d = $word.ActiveDocument

$d.Content.Find.Execute("52×15=780", $true, $true, $false, $false, $false, $true, 1, $false, "25×38=950", 2) | Out-Null
$d.Content.Find.Execute("94×11=1034", $true, $true, $false, $false, $false, $true, 1, $false, "24×21=504", 2) | Out-Null
$d.Content.Find.Execute("99×39=3861", $true, $true, $false, $false, $false, $true, 1, $false, "87×72=6264", 2) | Out-Null
$d.Content.Find.Execute("37×32=1184", $true, $true, $false, $false, $false, $true, 1, $false, "22×51=1122", 2) | Out-Null
$d.Content.Find.Execute("45×76=3420", $true, $true, $false, $false, $false, $true, 1, $false, "94×71=6674", 2) | Out-Null
$d.Content.Find.Execute("74×78=5772", $true, $true, $false, $false, $false, $true, 1, $false, "88×74=6512", 2) | Out-Null
$d.Content.Find.Execute("69×70=4830", $true, $true, $false, $false, $false, $true, 1, $false, "92×74=6808", 2) | Out-Null
$d.Content.Find.Execute("64×25=1600", $true, $true, $false, $false, $false, $true, 1, $false, "27×58=1566", 2) | Out-Null
$d.Content.Find.Execute("30×96=2880", $true, $true, $false, $false, $false, $true, 1, $false, "89×17=1513", 2) | Out-Null
$d.Content.Find.Execute("59×45=2655", $true, $true, $false, $false, $false, $true, 1, $false, "51×24=1224", 2) | Out-Null
$d.Content.Find.Execute("84×52=4368", $true, $true, $false, $false, $false, $true, 1, $false, "67×85=5695", 2) | Out-Null
$d.Content.Find.Execute("96×87=8352", $true, $true, $false, $false, $false, $true, 1, $false, "47×94=4418", 2) | Out-Null
$d.Content.Find.Execute("70×63=4410", $true, $true, $false, $false, $false, $true, 1, $false, "37×82=3034", 2) | Out-Null
$d.Content.Find.Execute("62×79=4898", $true, $true, $false, $false, $false, $true, 1, $false, "26×44=1144", 2) | Out-Null
$d.Content.Find.Execute("88×78=6864", $true, $true, $false, $false, $false, $true, 1, $false, "22×35=770", 2) | Out-Null
$d.Content.Find.Execute("28×85=2380", $true, $true, $false, $false, $false, $true, 1, $false, "56×63=3528", 2) | Out-Null
$d.Content.Find.Execute("88×83=7304", $true, $true, $false, $false, $false, $true, 1, $false, "39×19=741", 2) | Out-Null
$d.Content.Find.Execute("49×97=4753", $true, $true, $false, $false, $false, $true, 1, $false, "53×89=4717", 2) | Out-Null
$d.Content.Find.Execute("90×97=8730", $true, $true, $false, $false, $false, $true, 1, $false, "92×93=8556", 2) | Out-Null
$d.Content.Find.Execute("81×69=5589", $true, $true, $false, $false, $false, $true, 1, $false, "46×87=4002", 2) | Out-Null
$d.Content.Find.Execute("14×63=882", $true, $true, $false, $false, $false, $true, 1, $false, "58×96=5568", 2) | Out-Null
$d.Content.Find.Execute("71×32=2272", $true, $true, $false, $false, $false, $true, 1, $false, "12×50=600", 2) | Out-Null
$d.Content.Find.Execute("16×28=448", $true, $true, $false, $false, $false, $true, 1, $false, "41×57=2337", 2) | Out-Null
$d.Content.Find.Execute("52×31=1612", $true, $true, $false, $false, $false, $true, 1, $false, "66×95=6270", 2) | Out-Null
$d.Content.Find.Execute("64×54=3456", $true, $true, $false, $false, $false, $true, 1, $false, "78×68=5304", 2) | Out-Null
